$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "pv_pooled"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "pv_1,pv_2,pv_3,pv_4,pv_5"
$ws.Range("B3").Style = "Normal"

$ws.Range("A4").Value = "pvkat_pooled"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "pvkat_1,pvkat_2,pvkat_3,pvkat_4,pvkat_5"
$ws.Range("B4").Style = "Normal"
